{"js": "// Update the worksheet date and all twenty-five two-digit by two-digit\n// multiplication problems to the next day's generated values.\nconst replacements = [\n  [\"2024-10-03 Thursday\", \"2024-10-04 Friday\"],\n  [\"15\u00d764=\", \"96\u00d726=\"],\n  [\"84\u00d734=\", \"66\u00d739=\"],\n  [\"19\u00d776=\", \"31\u00d716=\"],\n  [\"60\u00d791=\", \"81\u00d794=\"],\n  [\"47\u00d721=\", \"77\u00d763=\"],\n  [\"70\u00d751=\", \"58\u00d773=\"],\n  [\"21\u00d776=\", \"63\u00d765=\"],\n  [\"66\u00d780=\", \"41\u00d795=\"],\n  [\"56\u00d791=\", \"65\u00d732=\"],\n  [\"25\u00d739=\", \"31\u00d781=\"],\n  [\"81\u00d756=\", \"30\u00d780=\"],\n  [\"52\u00d778=\", \"15\u00d791=\"],\n  [\"87\u00d747=\", \"76\u00d745=\"],\n  [\"36\u00d763=\", \"64\u00d742=\"],\n  [\"50\u00d759=\", \"61\u00d751=\"],\n  [\"40\u00d721=\", \"22\u00d759=\"],\n  [\"22\u00d752=\", \"42\u00d792=\"],\n  [\"23\u00d790=\", \"84\u00d755=\"],\n  [\"91\u00d716=\", \"47\u00d717=\"],\n  [\"16\u00d756=\", \"50\u00d792=\"],\n  [\"34\u00d789=\", \"29\u00d782=\"],\n  [\"60\u00d786=\", \"53\u00d763=\"],\n  [\"78\u00d777=\", \"50\u00d723=\"],\n  [\"71\u00d754=\", \"37\u00d797=\"],\n  [\"64\u00d736=\", \"17\u00d788=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all twenty-five two-digit by two-digit\n# multiplication problems to the next day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-03 Thursday\", \"2024-10-04 Friday\"),\n    @(\"15\u00d764=\", \"96\u00d726=\"),\n    @(\"84\u00d734=\", \"66\u00d739=\"),\n    @(\"19\u00d776=\", \"31\u00d716=\"),\n    @(\"60\u00d791=\", \"81\u00d794=\"),\n    @(\"47\u00d721=\", \"77\u00d763=\"),\n    @(\"70\u00d751=\", \"58\u00d773=\"),\n    @(\"21\u00d776=\", \"63\u00d765=\"),\n    @(\"66\u00d780=\", \"41\u00d795=\"),\n    @(\"56\u00d791=\", \"65\u00d732=\"),\n    @(\"25\u00d739=\", \"31\u00d781=\"),\n    @(\"81\u00d756=\", \"30\u00d780=\"),\n    @(\"52\u00d778=\", \"15\u00d791=\"),\n    @(\"87\u00d747=\", \"76\u00d745=\"),\n    @(\"36\u00d763=\", \"64\u00d742=\"),\n    @(\"50\u00d759=\", \"61\u00d751=\"),\n    @(\"40\u00d721=\", \"22\u00d759=\"),\n    @(\"22\u00d752=\", \"42\u00d792=\"),\n    @(\"23\u00d790=\", \"84\u00d755=\"),\n    @(\"91\u00d716=\", \"47\u00d717=\"),\n    @(\"16\u00d756=\", \"50\u00d792=\"),\n    @(\"34\u00d789=\", \"29\u00d782=\"),\n    @(\"60\u00d786=\", \"53\u00d763=\"),\n    @(\"78\u00d777=\", \"50\u00d723=\"),\n    @(\"71\u00d754=\", \"37\u00d797=\"),\n    @(\"64\u00d736=\", \"17\u00d788=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
